$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format while writing so that numeric-looking
# strings (e.g. "1.001", "0.9995", "4.061") are preserved exactly as text,
# matching the original inline-string cell contents instead of being
# auto-converted into floating point numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '24.481.87'
$ws.Range("E2").Value = '  -1.50%  '
$ws.Range("D3").Value = '1.652.21'
$ws.Range("E3").Value = '  -3.24%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").Value = '313.02'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").Value = '0.9995'
$ws.Range("D7").Value = '0.3662'
$ws.Range("E7").Value = '  -2.38%  '
$ws.Range("D8").Value = '46.71'
$ws.Range("E8").Value = '  -5.56%  '
$ws.Range("D9").Value = '0.3256'
$ws.Range("E9").Value = '  -5.69%  '
$ws.Range("D10").Value = '1.127'
$ws.Range("E10").Value = '  -7.37%  '
$ws.Range("D11").Value = '0.07054'
$ws.Range("E11").Value = '  -6.59%  '
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").Value = '5.978'
$ws.Range("E13").Value = '  -5.40%  '
$ws.Range("D14").Value = '19.43'
$ws.Range("E14").Value = '  -8.80%  '
$ws.Range("D15").Value = '6.625'
$ws.Range("E15").Value = '  -6.42%  '
$ws.Range("D16").Value = '1.655.20'
$ws.Range("E16").Value = '  -3.17%  '
$ws.Range("D17").Value = '0.00001044'
$ws.Range("E17").Value = '  -7.86%  '
$ws.Range("D18").Value = '0.06584'
$ws.Range("E18").Value = '  -2.11%  '
$ws.Range("D19").Value = '0.9987'
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").Value = '78.92'
$ws.Range("E20").Value = '  -6.79%  '
$ws.Range("D21").Value = '5.940'
$ws.Range("E21").Value = '  -7.28%  '
$ws.Range("D22").Value = '15.70'
$ws.Range("E22").Value = '  -9.36%  '
$ws.Range("D23").Value = '12.59'
$ws.Range("E23").Value = '  -4.32%  '
$ws.Range("D24").Value = '24.442.40'
$ws.Range("E24").Value = '  -1.73%  '
$ws.Range("D25").Value = '2.469'
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("D26").Value = '2.333'
$ws.Range("E26").Value = '  -16.44%  '
$ws.Range("D27").Value = '146.87'
$ws.Range("E27").Value = '  -2.47%  '
$ws.Range("D28").Value = '18.64'
$ws.Range("E28").Value = '  -8.89%  '
$ws.Range("D29").Value = '1.843.98'
$ws.Range("E29").Value = '  -2.94%  '
$ws.Range("D30").Value = '1.205'
$ws.Range("E30").Value = '  -3.37%  '
$ws.Range("D31").Value = '124.36'
$ws.Range("E31").Value = '  -6.41%  '
$ws.Range("D32").Value = '4.061'
$ws.Range("E32").Value = '  -4.23%  '
$ws.Range("D33").Value = '5.740'
$ws.Range("E33").Value = '  -16.29%  '
$ws.Range("D34").Value = '0.08466'
$ws.Range("E34").Value = '  -4.21%  '
$ws.Range("D35").Value = '1.654'
$ws.Range("E35").Value = '  -6.46%  '
$ws.Range("D36").Value = '12.17'
$ws.Range("E36").Value = '  -12.00%  '
$ws.Range("D37").Value = '5.231'
$ws.Range("E37").Value = '  -7.67%  '
$ws.Range("D38").Value = '1.267'
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.06037'
$ws.Range("E39").Value = '  -9.45%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.02239'
$ws.Range("E40").Value = '  -7.07%  '
$ws.Range("D41").Value = '0.2071'
$ws.Range("E41").Value = '  -7.40%  '
$ws.Range("D42").Value = '8.149'
$ws.Range("E42").Value = '  -12.31%  '
$ws.Range("D43").Value = '0.9993'
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("D44").Value = '0.5926'
$ws.Range("E44").Value = '  -8.30%  '
$ws.Range("D45").Value = '3.789'
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("D46").Value = '12.64'
$ws.Range("E46").Value = '  -9.02%  '
$ws.Range("D47").Value = '0.5632'
$ws.Range("E47").Value = '  -8.56%  '
$ws.Range("D48").Value = '123.13'
$ws.Range("E48").Value = '  -5.37%  '
$ws.Range("D49").Value = '1.953'
$ws.Range("E49").Value = '  -8.55%  '
$ws.Range("D50").Value = '0.06912'
$ws.Range("E50").Value = '  -5.41%  '
$ws.Range("D51").Value = '1.186'
$ws.Range("E51").Value = '  -2.79%  '

# Restore the default (Normal) style on the price column so no stray
# number-format styling is left behind on the cells.
$priceRange.Style = "Normal"
